$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44322
$ws.Range("D4").Value = 44327
$ws.Range("M4").Value = 60
$ws.Range("D5").Value = 44306
$ws.Range("D6").Value = 44302
$ws.Range("M6").Value = 80
$ws.Range("D7").Value = 44316
$ws.Range("M7").Value = 120
$ws.Range("D8").Value = 44330
$ws.Range("M8").Value = 60
$ws.Range("D9").Value = 44313
$ws.Range("M9").Value = 120
$ws.Range("D10").Value = 44323
$ws.Range("M10").Value = 80
